# Apply "repull data, push all data, mean calculation" edit:
# Update column F (dSF) values for specific rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -4
    8  = -4
    9  = -2
    12 = 0
    17 = 2
    18 = -6
    20 = -1
    21 = -3
    27 = -3
    28 = 3
    31 = 0
    32 = -2
    34 = 0
    36 = -2
    37 = 1
    38 = 4
    50 = 5
    54 = 2
    55 = 1
    57 = 0
    61 = -1
    70 = -4
    72 = -2
    73 = -1
    75 = 1
    76 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
